# Scen_ELC_CCS.xlsx - "Adjust naming convention and clean-up"
#  - rename the EPP* process-set wildcards to the P-TH* naming convention
#    used after merging the ELC SubRES files into a single workbook
#  - tidy up the sheet view (zoom back out, reset scroll position/selection)
#  - widen column B so the UC_N labels are fully visible

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UC_CCS")
$ws.Activate()

# --- Naming-convention rename: EPP*...* -> P-TH*...* -----------------------

# FLO_EMIS block (rows 3-8): the CCS process-set wildcard
$ws.Range("G3:G8").Value = "P-TH*CCS*"

# UC block (rows 16-20): per-plant process-set wildcards
$ws.Range("C16").Value = "P-TH*Moneypoint*"
$ws.Range("C17").Value = "P-TH*Cork*"
$ws.Range("C18").Value = "P-TH*Dublin*"
$ws.Range("C19").Value = "P-TH*Offaly*"
$ws.Range("C20").Value = "P-TH*Kilroot*"

# --- Cosmetic clean-up -------------------------------------------------

# Column B holds the long UC_N names (e.g. CCS_Moneypoint_MaxCap) - widen it
$ws.Columns.Item(2).ColumnWidth = 20.92

# Reset the view: zoom out to 70%, scroll back to the top and select the
# refreshed wildcard column
$excel.ActiveWindow.Zoom = 70
$ws.Range("C16:C20").Select()
